$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title subtitle EQUIPO 2 (row 2)
$ws.Range("D2").Value = "EQUIPO 2"

# Row 15
$ws.Range("B15").Value = "CP010"
$ws.Range("C15").Value = "Validar rechazo campo para publicar en el blog."
$ws.Range("D15").Value = "Comprobar que el usuario no pueda publicar contenido con emojis."
$ws.Range("E15").Value = "Ninguna"
$ws.Range("F15").Value = "1. Ingresar al blog.`n2. Click en el espacio para publicar contenido.`n3. Ingresar emojis.`n"
$ws.Range("G15").Value = "- Emojis ingresados.`n"
$ws.Range("H15").Value = "Ninguna"
$ws.Range("I15").Value = "- Mensaje de impedimento por ingreso datos no permitido"
$ws.Range("J15").Value = "Exitoso, el sistema no permite ingresar datos no validos"

# Row 16
$ws.Range("B16").Value = "CP011"
$ws.Range("C16").Value = "Validar publicación de enlaces web."
$ws.Range("D16").Value = "Comprobar la publicación de enlaces web en el blog"
$ws.Range("E16").Value = "Ninguna"
$ws.Range("F16").Value = "1. Escribir una dirección web de dominio `".com`", `".net`"  que se desea publicar.`n2. Copiar y pegar la dirección web de dominio `".com`", `".net`"  que se desea publicar.`n"
$ws.Range("G16").Value = "- Dirección web de dominio `".com`",  y `".net`"."
$ws.Range("H16").Value = "Ninguna"
$ws.Range("I16").Value = "- Dirección web ingresada"
$ws.Range("J16").Value = "Exitoso, el sistema permite realizar la publicación."

# Row 17
$ws.Range("B17").Value = "CP012"
$ws.Range("C17").Value = "Validar rechazo de publicación de enlaces web"
$ws.Range("D17").Value = "Comprobar que no se realice la publicación de enlaces no permitidos"
$ws.Range("E17").Value = "Ninguna"
$ws.Range("F17").Value = "1. Escribir una dirección web de dominio `".gov`" u otro dominio diferente a `".com`" o `".net`" que se desea publicar.`n2. Copiar y pegar la dirección web  de otro dominio diferente a `".com`" o `".net`" que se desea publicar.`n"
$ws.Range("G17").Value = "- Dirección web de dominio `".gov`" u otro dominio diferente a `".com`" o `".net`"."
$ws.Range("H17").Value = "Ninguna"
$ws.Range("I17").Value = "- Mensaje de impedimento de publicación por dominio o contenido no admitido."
$ws.Range("J17").Value = "Exitoso, el sistema no permite realizar la publicación."

# Row 18
$ws.Range("B18").Value = "CP013"
$ws.Range("C18").Value = "Validar edición de estilo y estrutura de texto en la entrada del blog."
$ws.Range("D18").Value = "Comprobar que el usuario pueda darle formato al texto que desea publicar"
$ws.Range("E18").Value = "Ninguna"
$ws.Range("F18").Value = "1. Ingresar un texto.`n2. Aplicar negrita al texto.`n3. Aplicar cursiva al texto.`n4. Aplicar subrayado al texto.`n5. Cambiar color al texto.`n6. Aumentar y disminuir tamaño al texto. `n7. Alinear texto (centro, derecha, izquierda)"
$ws.Range("G18").Value = "- Texto ingresado."
$ws.Range("H18").Value = "Ninguna"
$ws.Range("I18").Value = "- Texto en negrita.`n- Texto en cursiva.`n- Texto subrayado.`n- Texto a color.`n- Texto con más o menos tamaño.`n- Texto alineado."
$ws.Range("J18").Value = "Exitoso, el sistema permite publicar texto."

# Row 19
$ws.Range("B19").Value = "CP014"
$ws.Range("C19").Value = "Validar publicación del contenido de la entrada del Blog en redes sociales"
$ws.Range("D19").Value = "Comprobar que el sistema permita compartir una publicación en redes sociales."
$ws.Range("E19").Value = "Ninguna"
$ws.Range("F19").Value = "1. Elegir la publicación a compartir.`n2. Click en compartir a redes sociales.`n3. Click en compartir a través de Facebook.`n4. Ingresar los datos del usuario (Usuario y contraseña) y click en compartir.`n5. Click en compartir a través de Instagram.`n6. Ingresar los datos del usuario (Usuario y contraseña) y click en compartir.`n7. Click en compartir a través de Twitter.`n8. Ingresar los datos del usuario (Usuario y contraseña) y click en compartir. `n9.  Click en compartir a través de LinkedIn.`n10. Ingresar los datos del usuario (Usuario y contraseña) y click en compartir. "
$ws.Range("G19").Value = "- Publicación que se desea compartir."
$ws.Range("H19").Value = "Ninguna"
$ws.Range("I19").Value = "- Publicación a compartir."
$ws.Range("J19").Value = "Exitoso, el sistema permite compartir la publicación a redes sociales"

# Row 20
$ws.Range("B20").Value = "CP015"
$ws.Range("C20").Value = "Validar rechazo en publicación del contenido de la entrada del Blog en redes sociales"
$ws.Range("D20").Value = "Comprobar que el sistema permita no permita compartir una publicación en redes sociales sin el procedimiento correcto."
$ws.Range("E20").Value = "Ninguna"
$ws.Range("F20").Value = "1. Elegir la publicación a compartir.`n2. Click en compartir a redes sociales.`n3. Click en compartir a través de Facebook.`n4. No ingresar los datos del usuario (Usuario y contraseña) y click en compartir.`n5. Click en compartir a través de Instagram.`n6. No ingresar los datos del usuario (Usuario y contraseña) y click en compartir.`n7. Click en compartir a través de Twitter.`n8. No ingresar los datos del usuario (Usuario y contraseña) y click en compartir. `n9.  Click en compartir a través de LinkedIn.`n10. No ingresar los datos del usuario (Usuario y contraseña) y click en compartir. "
$ws.Range("G20").Value = "- Publicación que se desea compartir."
$ws.Range("H20").Value = "Ninguna"
$ws.Range("I20").Value = "- Mensaje de impedimento a compartir por no ingresar datos del usuario."
$ws.Range("J20").Value = "Exitoso, el sistema no permite compartir la publicación a redes sociales"

# Row 21
$ws.Range("B21").Value = "CP016"
$ws.Range("C21").Value = "Validar la impresión del contenido de la Entrada publicación del blog"
$ws.Range("D21").Value = "Comprobar que el sistema permita imprimir una publicación."
$ws.Range("E21").Value = "Ninguna"
$ws.Range("F21").Value = "1. Elegir una publicación con texto (incluídas direcciones web) o texto + imágenes para imprimir.`n2. Conectar la impresora con el ordenador.`n3. Click en imprimir.`n4. Seleccionar impresora.`n5. Aplicar formato deseado.`n6. Click en aceptar."
$ws.Range("G21").Value = "- Publicación que se desea imprimir."
$ws.Range("H21").Value = "Ninguna"
$ws.Range("I21").Value = "- Interfaz de impresión.`n- Vista previa de la impresión."
$ws.Range("J21").Value = "Exitoso, el sistema permite imprimir el contenido de una publicación."

# Row 22
$ws.Range("B22").Value = "CP017"
$ws.Range("C22").Value = "Validar rechazo de la impresión del contenido de la Entrada publicación del blog."
$ws.Range("D22").Value = "Comprobar que el sistema no permita imprimir publicaciones con videos."
$ws.Range("E22").Value = "Ninguna"
$ws.Range("F22").Value = "1. Elegir una publicación con uno o dos videos para imprimir.`n2. Conectar la impresora con el ordenador.`n3. Click en imprimir.`n4. Seleccionar impresora.`n5. Aplicar formato deseado.`n6. Click en aceptar."
$ws.Range("G22").Value = "- Publicación que se desea imprimir."
$ws.Range("H22").Value = "Ninguna"
$ws.Range("I22").Value = "- Mensaje de impedimento indicando que el contenido no se puede imprimir."
$ws.Range("J22").Value = "Exitoso, el sistema no permite imprimir un video."

# Row 23
$ws.Range("B23").Value = "CP018"
$ws.Range("C23").Value = "Validar el envio de una publicación por correo eletrónico."
$ws.Range("D23").Value = "Comprobar que el sistema permita enviar una publicación por correo electrónico "
$ws.Range("E23").Value = "Ninguna"
$ws.Range("F23").Value = "1. Elegir la publicación a enviar por correo.`n2. Click en `"Enviar vía correo electrónico`".`n3. Ingresar correo electrónico del Remitente.`n4. Ingresar correo electrónico del Destinatario.`n5. Ingresar Asunto (este campo es opcional)`n6. Click en enviar."
$ws.Range("G23").Value = "- Publicación a enviar por correo electrónico."
$ws.Range("H23").Value = "Ninguna"
$ws.Range("I23").Value = "- Vista previa del contenido de la publicación a enviar."
$ws.Range("J23").Value = "Exitoso, el sistema permite el envio de una publicación vía correo electrónico."

# Row 24
$ws.Range("B24").Value = "CP019"
$ws.Range("C24").Value = "Validar el rechazo de  enviar una publicación por correo eletrónico."
$ws.Range("D24").Value = "Comprobar que el sistema No permita enviar una publicación por correo electrónico sin el procedimiento correcto."
$ws.Range("E24").Value = "Ninguna"
$ws.Range("F24").Value = "1. Elegir la publicación a enviar por correo.`n2. Click en `"Enviar vía correo electrónico`".`n3. No ingresar correo electrónico del Remitente.`n4. No Ingresar correo electrónico del Destinatario.`n5. Ingresar Asunto (este campo es opcional)`n6. Click en enviar."
$ws.Range("G24").Value = "- Publicación a enviar por correo electrónico."
$ws.Range("H24").Value = "Ninguna"
$ws.Range("I24").Value = "- Mensaje de impedimento por falta de requerimientos obligatorios."
$ws.Range("J24").Value = "Exitoso, el sistema no permite el envio de la publicación vía correo electrónico.."

# Fix F13 text: "Cargar videos" -> "Insertar videos" (added last so its new shared string
# lands at the end of the shared-strings table, matching the authored edit order)
$ws.Range("F13").Value = "1. Seleccionar `"Insertar videos`".`n2. Seleccionar 2 videos en formato mp4 con tamaño mayor a 100Mb en algun video.`n3. Seleccionar 2 videos en un formato diferente a mp4 en al menos un video.`n4. Seleccionar 3 o más videos en cualquier formato."
